# slaten_justin.xlsx: regenerate save_data to use K (strikeouts) instead of
# the old "Strike#" source, and rewrite the recomputed std/mean-derived
# K values (column G) for each outing row.
#
# Column layout (row 1 header): A=idx, B=date, C=TB, D=PC, E=dS0, F=dSF,
# G=K, H=IP, I=I0, J=IF. Only column G (K) values change; every other
# column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 3
    22 = 3
    23 = 0
    24 = 3
    25 = 1
    26 = 1
    27 = 1
    28 = 3
    29 = 3
    30 = 1
    31 = 3
    32 = 2
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 5
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 2
    48 = 1
    49 = 1
    50 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
